$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("combinedrepos")

# Map of row number -> prefix to prepend to the existing flaw_rule value in column A.
$prefixes = @{
    2  = "B-"
    3  = "B-"
    4  = "B-"
    5  = "B-"
    6  = "V-"
    7  = "B-"
    8  = "B-"
    9  = "B-"
    10 = "V-"
    11 = "B-"
    12 = "B-"
    13 = "B-"
    14 = "B-"
    15 = "B-"
    16 = "B-"
    17 = "V-"
    18 = "B-"
    19 = "B-"
    20 = "B-"
    21 = "B-"
    22 = "B-"
    23 = "B-"
    24 = "B-"
    25 = "B-"
    26 = "B-"
    27 = "B-"
    28 = "V-"
    29 = "B-"
    30 = "B-"
    31 = "B-"
    32 = "B-"
    33 = "V-"
    34 = "B-"
    35 = "B-"
    36 = "B-"
    37 = "V-"
    38 = "B-"
    39 = "B-"
    40 = "B-"
    41 = "B-"
    42 = "B-"
    43 = "V-"
    44 = "B-"
    45 = "B-"
    46 = "V-"
}

foreach ($row in $prefixes.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Value2
    $cell.Value = $prefixes[$row] + $current
}
